$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "NormalString"
$ws.Range("B2").Value = "a"
$ws.Range("B3").Value = "b"
$ws.Range("B4").Value = "c"
$ws.Range("B5").Value = "d"

$ws.Columns.Item(2).ColumnWidth = 16.142857142857142

[void]$ws.Range("B6").Select()
